# "Fruta / hortaliza, semanal" update
#
# A new weekly price-report row is inserted into the price list at row 127
# (pushing the former rows 127-186 down to 128-187), for:
#   Vega Central Mapocho de Santiago / Ciruela / Angeleno / Primera
#   Fecha 2023-03-16 (serial 45001), 50 units, 12000/12000/12000,
#   $/bandeja 18 kilos granel, Región Metropolitana, 667 $/Kg, 18 Kg/unidad.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 127; everything below shifts down by one.
$ws.Rows("127:127").Insert()

$ws.Range("A127").Value = 9
$ws.Range("B127").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C127").Value = "Metropolitana"
$ws.Range("D127").Value = 45001
$ws.Range("E127").Value = 13
$ws.Range("F127").Value = "Fruta"
$ws.Range("G127").Value = 100103
$ws.Range("H127").Value = "Frutos de hueso (carozo)"
$ws.Range("I127").Value = 100103002
$ws.Range("J127").Value = "Ciruela"
$ws.Range("K127").Value = "Angeleno"
$ws.Range("L127").Value = "Primera"
$ws.Range("M127").Value = 50
$ws.Range("N127").Value = 12000
$ws.Range("O127").Value = 12000
$ws.Range("P127").Value = 12000
$ws.Range("Q127").Value = "$/bandeja 18 kilos granel"
$ws.Range("R127").Value = "Región Metropolitana"
$ws.Range("S127").Value = 667
$ws.Range("T127").Value = 18
